$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Highlight row 88 (A88:U88) with a solid red fill ---
$ws.Range("A88:U88").Interior.Color = 255

# --- Add note in U88 about mis-recorded approx fsr ---
$ws.Range("U88").Value = "actually have reason to think I recorded approx fsr wrong on this one, so need to redo. Because when putting in for next run saw it to still be 1.7."

# --- New data run rows for sg_rr_100_027 2023-12-08 17-44-55 ---

# Row 89: first attempt, distance = 1E-3
$ws.Range("A89").Value = "sg_rr_100_027 2023-12-08 17-44-55.csv"
$ws.Range("B89").Value = 0.01
$ws.Range("C89").Value = 1000
$ws.Range("D89").Value = 5001
$ws.Range("E89").Value = 1530
$ws.Range("F89").Value = 1570
$ws.Range("G89").Value = 0.001
$ws.Range("H89").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I89").Value = 1
$ws.Range("U89").Value = "seemed to find one peak in what looked like noise so increased prominence"

# Row 90: final run used for results, distance = 1.5E-3
$ws.Range("A90").Value = "sg_rr_100_027 2023-12-08 17-44-55.csv"
$ws.Range("B90").Value = 0.01
$ws.Range("C90").Value = 1000
$ws.Range("D90").Value = 5001
$ws.Range("E90").Value = 1530
$ws.Range("F90").Value = 1570
$ws.Range("G90").Value = 0.0015
$ws.Range("H90").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I90").Value = 1
$ws.Range("J90").Value = 0.98153846153846003
$ws.Range("K90").Value = 0.0043858818636388196
$ws.Range("L90").Value = "yes"
$ws.Range("M90").Value = 0.133427659342539
$ws.Range("N90").Value = 0.0039535320072782501
$ws.Range("O90").Value = 11915.510338869801
$ws.Range("P90").Value = 264.39216040517698
$ws.Range("Q90").Value = 169175684.93642601
$ws.Range("R90").Value = 11266552.773672201
$ws.Range("S90").Value = 100
$ws.Range("T90").Value = 0.1

# --- Update view state to match where the user ended up scrolled/selected ---
$ws.Range("A61").Select()
